# Update the "payload" column (B) for rows 15-27 from 0 to 1,
# then move the active selection to N11 to mirror the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("slipstream_event_schedule")

$ws.Range("B15:B27").Value = 1

$ws.Range("N11").Select()
